# Add row 7: a new Id/Question/Answer entry for item "6", mirroring the
# pattern used by rows 2-6 (timestamp+question id / "N question" / "N").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "2023-10-17 13:42:35 6 question_6_2"
$ws.Range("B7").Value = "6 question"

# Column C holds the answer "6" as TEXT (matching existing cells C2:C6,
# which all store their digit as a shared string, not a number). Typing
# "6" straight into Value would be auto-coerced to a numeric cell, so we
# stage it in a scratch cell formatted as Text, then copy only the VALUE
# (PasteSpecial xlPasteValues = -4163) into C7 - that carries the text
# type over without carrying the scratch cell's number format/style.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "6"
$ws.Range("Z1").Copy()
$ws.Range("C7").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
